# Regenerate the Handback report: the handoff run that produced this
# workbook was re-run, producing fresh file GUIDs, a fresh content hash,
# and fresh timestamps. Propagate the new identifiers/timestamps to every
# sheet (Overview, zh-cn, de-de) and to the hyperlink display text that
# mirrors them.

$wb = $excel.ActiveWorkbook

$newGuid1 = "7a6ff401-3d4c-4db2-a3e2-a2ff2732b77b"
$newGuid2 = "ffffa549b742-c9a4-44e5-90d1-e0b93a6912c2"

$newHash = "36333619e9ac0be9a3ef7c6d3f3d1fb72be4b193"

$newRowTime       = "2016-08-19 17:06:32"
$newHandoffTime   = "2016-08-19 17:06:28"
$newHandbackZh    = "2016-08-19 17:06:44"
$newHandbackDe    = "2016-08-19 17:06:51"

$newFile1 = $newGuid1 + ".md"
$newFile2 = $newGuid2 + ".md"
$newPath1 = "e2e\" + $newGuid1 + ".md"
$newPath2 = "e2e\" + $newGuid2 + ".md"

$newXlfZh1 = $newGuid1 + "." + $newHash + ".zh-cn.xlf"
$newXlfDe = $newGuid1 + "." + $newHash + ".de-de.xlf"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $newFile1
$wsOverview.Range("B2").Value = $newPath1
$wsOverview.Range("G2").Value = $newRowTime

$wsOverview.Range("A3").Value = $newFile2
$wsOverview.Range("B3").Value = $newPath2
$wsOverview.Range("G3").Value = $newRowTime

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = $newPath1
    }
    if ($addr -eq '$B$3') {
        $hl.TextToDisplay = $newPath2
    }
}

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = $newFile1
$wsZh.Range("G2").Value = $newXlfZh1
$wsZh.Range("H2").Value = $newHandoffTime
$wsZh.Range("I2").Value = $newFile1
$wsZh.Range("J2").Value = $newXlfZh1
$wsZh.Range("K2").Value = $newHandbackZh

$wsZh.Range("A3").Value = $newFile2
$wsZh.Range("G3").Value = $newXlfZh1
$wsZh.Range("H3").Value = $newHandoffTime
$wsZh.Range("I3").Value = $newFile2
$wsZh.Range("J3").Value = $newXlfZh1
$wsZh.Range("K3").Value = $newHandbackZh

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    }
    if ($addr -eq '$I$2') {
        $hl.TextToDisplay = $newFile1
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    }
    if ($addr -eq '$I$3') {
        $hl.TextToDisplay = $newFile2
    }
}

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = $newFile1
$wsDe.Range("G2").Value = $newXlfDe
$wsDe.Range("H2").Value = $newRowTime
$wsDe.Range("I2").Value = $newFile1
$wsDe.Range("J2").Value = $newXlfDe
$wsDe.Range("K2").Value = $newHandbackDe

$wsDe.Range("A3").Value = $newFile2
$wsDe.Range("G3").Value = $newXlfDe
$wsDe.Range("H3").Value = $newRowTime
$wsDe.Range("I3").Value = $newFile2
$wsDe.Range("J3").Value = $newXlfDe
$wsDe.Range("K3").Value = $newHandbackDe

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newFile1
    }
    if ($addr -eq '$I$2') {
        $hl.TextToDisplay = $newFile1
    }
    if ($addr -eq '$A$3') {
        $hl.TextToDisplay = $newFile2
    }
    if ($addr -eq '$I$3') {
        $hl.TextToDisplay = $newFile2
    }
}
